{"js": "// Update the 25 division problems in the practice table to the newly\n// generated set of operands, leaving the date line and everything else\n// untouched. Each \"before\" text is unique in the document, so a plain\n// (non-wildcard) search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"456\u00f74=\", \"552\u00f75=\"],\n  [\"987\u00f75=\", \"567\u00f74=\"],\n  [\"999\u00f76=\", \"403\u00f76=\"],\n  [\"450\u00f72=\", \"649\u00f78=\"],\n  [\"634\u00f74=\", \"150\u00f74=\"],\n  [\"246\u00f75=\", \"810\u00f77=\"],\n  [\"142\u00f78=\", \"166\u00f77=\"],\n  [\"915\u00f75=\", \"487\u00f75=\"],\n  [\"460\u00f74=\", \"809\u00f77=\"],\n  [\"457\u00f75=\", \"602\u00f72=\"],\n  [\"199\u00f76=\", \"551\u00f76=\"],\n  [\"977\u00f72=\", \"352\u00f79=\"],\n  [\"493\u00f75=\", \"927\u00f76=\"],\n  [\"703\u00f73=\", \"209\u00f73=\"],\n  [\"979\u00f75=\", \"186\u00f72=\"],\n  [\"622\u00f72=\", \"743\u00f73=\"],\n  [\"158\u00f77=\", \"881\u00f75=\"],\n  [\"438\u00f74=\", \"186\u00f74=\"],\n  [\"797\u00f78=\", \"892\u00f79=\"],\n  [\"565\u00f76=\", \"351\u00f74=\"],\n  [\"775\u00f74=\", \"898\u00f72=\"],\n  [\"776\u00f76=\", \"600\u00f73=\"],\n  [\"569\u00f72=\", \"118\u00f75=\"],\n  [\"723\u00f78=\", \"271\u00f78=\"],\n  [\"778\u00f75=\", \"213\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 division problems in the practice table to the newly\n# generated set of operands, leaving the date line and everything else\n# untouched. Each \"before\" text is unique in the document, so a plain\n# Find/Replace per pair (ReplaceAll) is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"456\u00f74=\", \"552\u00f75=\"),\n  @(\"987\u00f75=\", \"567\u00f74=\"),\n  @(\"999\u00f76=\", \"403\u00f76=\"),\n  @(\"450\u00f72=\", \"649\u00f78=\"),\n  @(\"634\u00f74=\", \"150\u00f74=\"),\n  @(\"246\u00f75=\", \"810\u00f77=\"),\n  @(\"142\u00f78=\", \"166\u00f77=\"),\n  @(\"915\u00f75=\", \"487\u00f75=\"),\n  @(\"460\u00f74=\", \"809\u00f77=\"),\n  @(\"457\u00f75=\", \"602\u00f72=\"),\n  @(\"199\u00f76=\", \"551\u00f76=\"),\n  @(\"977\u00f72=\", \"352\u00f79=\"),\n  @(\"493\u00f75=\", \"927\u00f76=\"),\n  @(\"703\u00f73=\", \"209\u00f73=\"),\n  @(\"979\u00f75=\", \"186\u00f72=\"),\n  @(\"622\u00f72=\", \"743\u00f73=\"),\n  @(\"158\u00f77=\", \"881\u00f75=\"),\n  @(\"438\u00f74=\", \"186\u00f74=\"),\n  @(\"797\u00f78=\", \"892\u00f79=\"),\n  @(\"565\u00f76=\", \"351\u00f74=\"),\n  @(\"775\u00f74=\", \"898\u00f72=\"),\n  @(\"776\u00f76=\", \"600\u00f73=\"),\n  @(\"569\u00f72=\", \"118\u00f75=\"),\n  @(\"723\u00f78=\", \"271\u00f78=\"),\n  @(\"778\u00f75=\", \"213\u00f78=\")\n)\n\n$wdReplaceAll = 2\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $wdReplaceAll)\n}\n"}
